# Insert a new data row at row 131 (pushing existing rows 131-161 down to 132-162)
# and populate it with the new record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 131, shifting rows 131..161 down to 132..162
$ws.Rows.Item(131).Insert()

# Populate the newly inserted row 131 with the new record's values
$ws.Cells.Item(131, 1).Value = 10
$ws.Cells.Item(131, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(131, 3).Value = "La Araucanía"
$ws.Cells.Item(131, 4).Value = 44798
$ws.Cells.Item(131, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(131, 5).Value = 9
$ws.Cells.Item(131, 6).Value = "Fruta"
$ws.Cells.Item(131, 7).Value = 100104
$ws.Cells.Item(131, 8).Value = "Frutos de pepita"
$ws.Cells.Item(131, 9).Value = 100104001
$ws.Cells.Item(131, 10).Value = "Granada"
$ws.Cells.Item(131, 11).Value = "Wonderfull"
$ws.Cells.Item(131, 12).Value = "Primera"
$ws.Cells.Item(131, 13).Value = 65
$ws.Cells.Item(131, 14).Value = 15000
$ws.Cells.Item(131, 15).Value = 15000
$ws.Cells.Item(131, 16).Value = 15000
$ws.Cells.Item(131, 17).Value = "`$/bandeja 10 kilos granel"
$ws.Cells.Item(131, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(131, 19).Value = 1500
$ws.Cells.Item(131, 20).Value = 10
